$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Reshape rows: the blank row 25 (between the "Apollo client cache" header
# and its single data row) is removed, and a new 12-row block is inserted
# after the data row so a brand-new "Apollo client - persisted queries"
# sub-table can be added, pushing the trailing summary rows down.
$ws.Rows("25:25").Delete()
$ws.Rows("32:43").Insert()

# --- Copy the formatting of the analogous existing blocks onto the new rows
# so fills/number formats match the rest of the sheet.
$ws.Range("A24:E24").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)
$ws.Range("A11:E11").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)
$ws.Range("A12:E12").Copy()
$ws.Range("A30:E30").PasteSpecial(-4122)
$ws.Range("B13:E13").Copy()
$ws.Range("B31:E31").PasteSpecial(-4122)
$ws.Range("A14:E14").Copy()
$ws.Range("A32:E32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 28: new sub-section header
$ws.Range("A28").Value = "Apollo client - persisted queries - without cdn - own server memory cache"

# --- Row 29: Frontend measurement
$ws.Range("A29").Value = "Frontend"
$ws.Range("B29").Value = 35
$ws.Range("C29").Formula = "=C30+1.5"
$ws.Range("D29:E29").Formula = "=D30+1.5"

# --- Row 30: Backend measurement
$ws.Range("A30").Value = "Backend"
$ws.Range("B30").Value = 33
$ws.Range("C30").Value = 181
$ws.Range("D30").Value = 35
$ws.Range("E30").Value = 55

# --- Row 31: averages of the two rows above
$ws.Range("B31").Formula = "=AVERAGE(B29:B30)"
$ws.Range("C31:E31").Formula = "=AVERAGE(C29:C30)"

# --- Row 32: totals row for the new block
$ws.Range("A32").Value = "mount tot inladen data"
$ws.Range("B32").Value = 56
$ws.Range("C32").Value = 216
$ws.Range("D32").Value = 77
$ws.Range("E32").Value = 102

# --- Match the author's final cursor position
$ws.Range("F38").Select() | Out-Null
